$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header A1 from "id_mapel" to "nomor"
$ws.Range("A1").Value = "nomor"

# Fill in the new row-numbering column (A2:A5) with sequential numbers
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# Highlight the header row with a yellow fill
$ws.Range("A1:D1").Interior.Color = 65535

# Update the visible scroll position / selection saved with the sheet
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("L25").Select()
